# MAJ timecodes et commentaires pour montage
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / value updates -------------------------------------------------

# New comment in G column, row 3
$ws.Range("G3").Value = "Les times codes sont pas fixes. Privilégier la synchronisation arrivée de la deuxième araignée avec la 2'' mélodie de la boite à musique"

# Updated total comment, row 4
$ws.Range("G4").Value = "(total es 3 précédents plans=504)"

# New comment about the pendulum tic-tac, row 9
$ws.Range("G9").Value = """/!\ à synchroniser le tic tac au mouvement du balancier"""

# Status text updates in F column
$ws.Range("F13").Value = "Rendu"
$ws.Range("F15").Value = "Rendu idéalement à refaire"
$ws.Range("F16").Value = "Rendu à faire"
$ws.Range("F17").Value = "Rendu"

# Shot-name corrections in A column
$ws.Range("A17").Value = "Cerf qui hume"
$ws.Range("A19").Value = "Contre plongée"

# New timecode comment, row 21
$ws.Range("G21").Value = "timecode fondu au noir: 2:22"

# --- Fill-colour (status) updates, reusing existing formats ---------------
# Donor cells F2 (yellow), F5 (light green) and F6 (dark green) keep their
# original colour throughout, so they can be used as format sources.

$ws.Range("F5").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("F17").PasteSpecial(-4122)

$ws.Range("F6").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("F16").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Column width & selection ---------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 22.14

$ws.Range("G20").Select()
